$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column G width (raw stored width target ~11.09 chars)
$ws.Columns.Item(7).ColumnWidth = 10.25

# Row 1: H1 gets a date value; copy date style from G1 first so it reuses the same style index
$ws.Cells.Item(1, 7).Copy() | Out-Null
$ws.Cells.Item(1, 8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 8).Value = 43952
$excel.CutCopyMode = 0

$ws.Cells.Item(2, 7).Value = 74966
$ws.Cells.Item(2, 8).Style = "Normal"
$ws.Cells.Item(2, 8).Value = 28183
$ws.Cells.Item(2, 9).Clear()
$ws.Cells.Item(3, 7).Value = 10313
$ws.Cells.Item(3, 8).Style = "Normal"
$ws.Cells.Item(3, 8).Value = 9404
$ws.Cells.Item(3, 9).Clear()
$ws.Cells.Item(4, 7).Value = 17671
$ws.Cells.Item(4, 8).Style = "Normal"
$ws.Cells.Item(4, 8).Value = 12436
$ws.Cells.Item(4, 9).Clear()
$ws.Cells.Item(5, 7).Value = 52581
$ws.Cells.Item(5, 8).Style = "Normal"
$ws.Cells.Item(5, 8).Value = 42909
$ws.Cells.Item(5, 9).Clear()
$ws.Cells.Item(6, 7).Value = 325343
$ws.Cells.Item(6, 8).Style = "Normal"
$ws.Cells.Item(6, 8).Value = 318064
$ws.Cells.Item(6, 9).Clear()
$ws.Cells.Item(7, 7).Value = 38662
$ws.Cells.Item(7, 8).Style = "Normal"
$ws.Cells.Item(7, 8).Value = 28322
$ws.Cells.Item(7, 9).Clear()
$ws.Cells.Item(8, 7).Value = 33041
$ws.Cells.Item(8, 8).Style = "Normal"
$ws.Cells.Item(8, 8).Value = 36166
$ws.Cells.Item(8, 9).Clear()
$ws.Cells.Item(9, 7).Value = 7947
$ws.Cells.Item(9, 8).Style = "Normal"
$ws.Cells.Item(9, 8).Value = 6183
$ws.Cells.Item(9, 9).Clear()
$ws.Cells.Item(10, 7).Value = 8708
$ws.Cells.Item(10, 8).Style = "Normal"
$ws.Cells.Item(10, 8).Value = 8133
$ws.Cells.Item(10, 9).Clear()
$ws.Cells.Item(11, 7).Value = 433103
$ws.Cells.Item(11, 8).Style = "Normal"
$ws.Cells.Item(11, 8).Value = 173191
$ws.Cells.Item(11, 9).Clear()
$ws.Cells.Item(12, 7).Value = 266565
$ws.Cells.Item(12, 8).Style = "Normal"
$ws.Cells.Item(12, 8).Value = 226884
$ws.Cells.Item(12, 9).Clear()
$ws.Cells.Item(13, 7).Value = 22495
$ws.Cells.Item(13, 8).Style = "Normal"
$ws.Cells.Item(13, 8).Value = 16112
$ws.Cells.Item(13, 9).Clear()
$ws.Cells.Item(14, 7).Value = 8827
$ws.Cells.Item(14, 8).Style = "Normal"
$ws.Cells.Item(14, 8).Value = 7194
$ws.Cells.Item(14, 9).Clear()
$ws.Cells.Item(15, 7).Value = 81596
$ws.Cells.Item(15, 8).Style = "Normal"
$ws.Cells.Item(15, 8).Value = 74476
$ws.Cells.Item(15, 9).Clear()
$ws.Cells.Item(16, 7).Value = 55774
$ws.Cells.Item(16, 8).Style = "Normal"
$ws.Cells.Item(16, 8).Value = 43777
$ws.Cells.Item(16, 9).Clear()
$ws.Cells.Item(17, 7).Value = 27220
$ws.Cells.Item(17, 8).Style = "Normal"
$ws.Cells.Item(17, 8).Value = 24693
$ws.Cells.Item(17, 9).Clear()
$ws.Cells.Item(18, 7).Value = 24483
$ws.Cells.Item(18, 8).Style = "Normal"
$ws.Cells.Item(18, 8).Value = 18281
$ws.Cells.Item(18, 9).Clear()
$ws.Cells.Item(19, 7).Value = 91223
$ws.Cells.Item(19, 8).Style = "Normal"
$ws.Cells.Item(19, 8).Value = 80060
$ws.Cells.Item(19, 9).Clear()
$ws.Cells.Item(20, 7).Value = 66141
$ws.Cells.Item(20, 8).Style = "Normal"
$ws.Cells.Item(20, 8).Value = 52137
$ws.Cells.Item(20, 9).Clear()
$ws.Cells.Item(21, 7).Value = 7661
$ws.Cells.Item(21, 8).Style = "Normal"
$ws.Cells.Item(21, 8).Value = 16175
$ws.Cells.Item(21, 9).Clear()
$ws.Cells.Item(22, 7).Value = 37925
$ws.Cells.Item(22, 8).Style = "Normal"
$ws.Cells.Item(22, 8).Value = 65262
$ws.Cells.Item(22, 9).Clear()
$ws.Cells.Item(23, 7).Value = 71358
$ws.Cells.Item(23, 8).Style = "Normal"
$ws.Cells.Item(23, 8).Value = 55448
$ws.Cells.Item(23, 9).Clear()
$ws.Cells.Item(24, 7).Value = 82004
$ws.Cells.Item(24, 8).Style = "Normal"
$ws.Cells.Item(24, 8).Value = 68952
$ws.Cells.Item(24, 9).Clear()
$ws.Cells.Item(25, 7).Value = 48595
$ws.Cells.Item(25, 8).Style = "Normal"
$ws.Cells.Item(25, 8).Value = 47134
$ws.Cells.Item(25, 9).Clear()
$ws.Cells.Item(26, 7).Value = 29906
$ws.Cells.Item(26, 8).Style = "Normal"
$ws.Cells.Item(26, 8).Value = 24810
$ws.Cells.Item(26, 9).Clear()
$ws.Cells.Item(27, 7).Value = 55299
$ws.Cells.Item(27, 8).Style = "Normal"
$ws.Cells.Item(27, 8).Value = 49402
$ws.Cells.Item(27, 9).Clear()
$ws.Cells.Item(28, 7).Value = 7052
$ws.Cells.Item(28, 8).Style = "Normal"
$ws.Cells.Item(28, 8).Value = 4263
$ws.Cells.Item(28, 9).Clear()
$ws.Cells.Item(29, 7).Value = 98941
$ws.Cells.Item(29, 8).Style = "Normal"
$ws.Cells.Item(29, 8).Value = 84716
$ws.Cells.Item(29, 9).Clear()
$ws.Cells.Item(30, 7).Value = 6274
$ws.Cells.Item(30, 8).Style = "Normal"
$ws.Cells.Item(30, 8).Value = 4689
$ws.Cells.Item(30, 9).Clear()
$ws.Cells.Item(31, 7).Value = 8229
$ws.Cells.Item(31, 8).Style = "Normal"
$ws.Cells.Item(31, 8).Value = 6555
$ws.Cells.Item(31, 9).Clear()
$ws.Cells.Item(32, 7).Value = 15001
$ws.Cells.Item(32, 8).Style = "Normal"
$ws.Cells.Item(32, 8).Value = 11834
$ws.Cells.Item(32, 9).Clear()
$ws.Cells.Item(33, 7).Value = 71966
$ws.Cells.Item(33, 8).Style = "Normal"
$ws.Cells.Item(33, 8).Value = 87540
$ws.Cells.Item(33, 9).Clear()
$ws.Cells.Item(34, 7).Value = 12093
$ws.Cells.Item(34, 8).Style = "Normal"
$ws.Cells.Item(34, 8).Value = 16801
$ws.Cells.Item(34, 9).Clear()
$ws.Cells.Item(35, 7).Value = 42541
$ws.Cells.Item(35, 8).Style = "Normal"
$ws.Cells.Item(35, 8).Value = 30735
$ws.Cells.Item(35, 9).Clear()
$ws.Cells.Item(36, 7).Value = 219413
$ws.Cells.Item(36, 8).Style = "Normal"
$ws.Cells.Item(36, 8).Value = 195242
$ws.Cells.Item(36, 9).Clear()
$ws.Cells.Item(37, 7).Value = 93599
$ws.Cells.Item(37, 8).Style = "Normal"
$ws.Cells.Item(37, 8).Value = 61046
$ws.Cells.Item(37, 9).Clear()
$ws.Cells.Item(38, 7).Value = 52500
$ws.Cells.Item(38, 8).Style = "Normal"
$ws.Cells.Item(38, 8).Value = 68237
$ws.Cells.Item(38, 9).Clear()
$ws.Cells.Item(39, 7).Value = 49300
$ws.Cells.Item(39, 8).Style = "Normal"
$ws.Cells.Item(39, 8).Value = 45102
$ws.Cells.Item(39, 9).Clear()
$ws.Cells.Item(40, 7).Value = 127896
$ws.Cells.Item(40, 8).Style = "Normal"
$ws.Cells.Item(40, 8).Value = 96603
$ws.Cells.Item(40, 9).Clear()
$ws.Cells.Item(41, 7).Value = 17286
$ws.Cells.Item(41, 8).Style = "Normal"
$ws.Cells.Item(41, 8).Value = 21673
$ws.Cells.Item(41, 9).Clear()
$ws.Cells.Item(42, 7).Value = 13084
$ws.Cells.Item(42, 8).Style = "Normal"
$ws.Cells.Item(42, 8).Value = 9109
$ws.Cells.Item(42, 9).Clear()
$ws.Cells.Item(43, 7).Value = 66438
$ws.Cells.Item(43, 8).Style = "Normal"
$ws.Cells.Item(43, 8).Value = 46747
$ws.Cells.Item(43, 9).Clear()
$ws.Cells.Item(44, 7).Value = 5535
$ws.Cells.Item(44, 8).Style = "Normal"
$ws.Cells.Item(44, 8).Value = 3756
$ws.Cells.Item(44, 9).Clear()
$ws.Cells.Item(45, 7).Value = 42805
$ws.Cells.Item(45, 8).Style = "Normal"
$ws.Cells.Item(45, 8).Value = 37319
$ws.Cells.Item(45, 9).Clear()
$ws.Cells.Item(46, 7).Value = 254084
$ws.Cells.Item(46, 8).Style = "Normal"
$ws.Cells.Item(46, 8).Value = 247179
$ws.Cells.Item(46, 9).Clear()
$ws.Cells.Item(47, 7).Value = 11738
$ws.Cells.Item(47, 8).Style = "Normal"
$ws.Cells.Item(47, 8).Value = 9057
$ws.Cells.Item(47, 9).Clear()
$ws.Cells.Item(48, 7).Value = 5117
$ws.Cells.Item(48, 8).Style = "Normal"
$ws.Cells.Item(48, 8).Value = 3759
$ws.Cells.Item(48, 9).Clear()
$ws.Cells.Item(49, 7).Value = 103
$ws.Cells.Item(49, 8).Style = "Normal"
$ws.Cells.Item(49, 8).Value = 11
$ws.Cells.Item(49, 9).Clear()
$ws.Cells.Item(50, 7).Value = 72488
$ws.Cells.Item(50, 8).Style = "Normal"
$ws.Cells.Item(50, 8).Value = 61138
$ws.Cells.Item(50, 9).Clear()
$ws.Cells.Item(51, 7).Value = 139505
$ws.Cells.Item(51, 8).Style = "Normal"
$ws.Cells.Item(51, 8).Value = 109167
$ws.Cells.Item(51, 9).Clear()
$ws.Cells.Item(52, 7).Value = 29818
$ws.Cells.Item(52, 8).Style = "Normal"
$ws.Cells.Item(52, 8).Value = 12996
$ws.Cells.Item(52, 9).Clear()
$ws.Cells.Item(53, 7).Value = 49993
$ws.Cells.Item(53, 8).Style = "Normal"
$ws.Cells.Item(53, 8).Value = 38002
$ws.Cells.Item(53, 9).Clear()
$ws.Cells.Item(54, 7).Value = 3497
$ws.Cells.Item(54, 8).Style = "Normal"
$ws.Cells.Item(54, 8).Value = 2026
$ws.Cells.Item(54, 9).Clear()

# Update selection to M10
$ws.Range("M10").Select()
